# Actualización desde MV -datos-
# Add the new "01-09-2021" row (row 24) to the demand-components table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A holds the period label as plain text (same as all the rows
# above it), even though it looks like a date. Typing the string directly
# via .Value would make Excel auto-convert it to a date serial number, so
# instead we build it as a text formula in a scratch cell and paste back
# only the *value* - this keeps the cell a plain shared string with no
# special number formatting/style, just like the existing rows.
$scratch = $ws.Range("K1")
$scratch.Formula = '="01-09-2021"'
$scratch.Copy()
$ws.Range("A24").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
$excel.CutCopyMode = $false

$ws.Range("D24").Value = 16.4
$ws.Range("E24").Value = 2.4
$ws.Range("H24").Value = 11.7
$ws.Range("I24").Value = 2.5
